$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows (1-based worksheet row numbers) where column B currently says "high risk"
# but should also be updated to "mixed or unspecified population".
$highRiskRowsToChange = @(26,50,198,202,328,354,402,481,545,585,649,654)

$newValue = "mixed or unspecified population"
$oldRegularValue = "regular, several popualtions, or unspecified"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $current = $cell.Value()
    if ($current -eq $oldRegularValue) {
        $cell.Value = $newValue
    } elseif ($highRiskRowsToChange -contains $r) {
        $cell.Value = $newValue
    }
}
